# Apply scheduled-runner price/profit updates across all job leve tables.
$wb = $excel.ActiveWorkbook

# ==== ALC sheet ====
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 380232.62
$ws.Range("J17").Value = 387307.75
$ws.Range("L17").Value = 1161923.25
$ws.Range("N17").Value = -1162259.25
# Row 19
$ws.Range("H19").Value = 1362.8
$ws.Range("I19").Value = 1328.5
$ws.Range("J19").Value = 1500
$ws.Range("K19").Value = 1328.5
$ws.Range("L19").Value = 1500
$ws.Range("M19").Value = -1153.5
$ws.Range("N19").Value = -1850
# Row 28
$ws.Range("H28").Value = 403799.8
$ws.Range("I28").Value = 669000
$ws.Range("K28").Value = 669000
$ws.Range("M28").Value = -668515
# Row 43
$ws.Range("H43").Value = 14793.823
$ws.Range("I43").Value = 11428.286
$ws.Range("J43").Value = 17149.7
$ws.Range("K43").Value = 11428.286
$ws.Range("L43").Value = 17149.7
$ws.Range("M43").Value = -11359.286
$ws.Range("N43").Value = -17287.7
# Row 51
$ws.Range("H51").Value = 3642.4546
$ws.Range("I51").Value = 4148.875
$ws.Range("J51").Value = 2292
$ws.Range("K51").Value = 4148.875
$ws.Range("L51").Value = 2292
$ws.Range("M51").Value = -3664.875
$ws.Range("N51").Value = -3260
# Row 107
$ws.Range("H107").Value = 1925.1852
$ws.Range("I107").Value = 2003.2
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 2003.2
$ws.Range("L107").Value = 950
$ws.Range("M107").Value = -83.20000000000005
$ws.Range("N107").Value = -4790
# Row 113
$ws.Range("H113").Value = 4496.9
$ws.Range("I113").Value = 4074.75
$ws.Range("J113").Value = 4778.3335
$ws.Range("K113").Value = 4074.75
$ws.Range("L113").Value = 4778.3335
$ws.Range("M113").Value = -820.75
$ws.Range("N113").Value = -11286.3335
# Row 121
$ws.Range("H121").Value = 1696.7273
$ws.Range("J121").Value = 1736.4
$ws.Range("L121").Value = 5209.200000000001
$ws.Range("N121").Value = -8703.200000000001
# Row 138
$ws.Range("H138").Value = 4200.79
$ws.Range("I138").Value = 2882.6667
$ws.Range("J138").Value = 4490.1343
$ws.Range("K138").Value = 8648.000100000001
$ws.Range("L138").Value = 13470.4029
$ws.Range("M138").Value = -3508.000100000001
$ws.Range("N138").Value = -23750.4029
# Row 141
$ws.Range("H141").Value = 4581.8335
$ws.Range("I141").Value = 4123.4287
$ws.Range("K141").Value = 12370.2861
$ws.Range("M141").Value = -7190.286100000001

# ==== ARM sheet ====
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13650946
$ws.Range("I32").Value = 13545562
$ws.Range("K32").Value = 13545562
$ws.Range("M32").Value = -13545275
# Row 45
$ws.Range("H45").Value = 5085.9
$ws.Range("I45").Value = 5053.5
$ws.Range("K45").Value = 5053.5
$ws.Range("M45").Value = -4676.5
# Row 61
$ws.Range("H61").Value = 2974.0244
$ws.Range("I61").Value = 2751.1562
$ws.Range("K61").Value = 2751.1562
$ws.Range("M61").Value = -2539.1562
# Row 74
$ws.Range("H74").Value = 2161.4
$ws.Range("I74").Value = 2151.8
$ws.Range("K74").Value = 2151.8
$ws.Range("M74").Value = -1277.8
# Row 77
$ws.Range("H77").Value = 2161.4
$ws.Range("I77").Value = 2151.8
$ws.Range("K77").Value = 10759
$ws.Range("M77").Value = -6391
# Row 92
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992
# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
# Row 122
$ws.Range("H122").Value = 2100.762
$ws.Range("I122").Value = 1840.8379
$ws.Range("K122").Value = 5522.5137
$ws.Range("M122").Value = -3072.5137
# Row 132
$ws.Range("H132").Value = 323990.75
$ws.Range("I132").Value = 418061.97
$ws.Range("J132").Value = 1460.8572
$ws.Range("K132").Value = 1254185.91
$ws.Range("L132").Value = 4382.571599999999
$ws.Range("M132").Value = -1251655.91
$ws.Range("N132").Value = -9442.571599999999
# Row 136
$ws.Range("H136").Value = 2974.0244
$ws.Range("I136").Value = 2751.1562
$ws.Range("K136").Value = 8253.4686
$ws.Range("M136").Value = -5703.4686

# ==== BSM sheet ====
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1094810
$ws.Range("I134").Value = 1192189.8
$ws.Range("K134").Value = 3576569.4
$ws.Range("M134").Value = -3574034.4

# ==== CRP sheet ====
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 1950.5
$ws.Range("I22").Value = 1126.5
$ws.Range("J22").Value = 5246.5
$ws.Range("K22").Value = 1126.5
$ws.Range("L22").Value = 5246.5
$ws.Range("M22").Value = -776.5
$ws.Range("N22").Value = -5946.5
# Row 31
$ws.Range("H31").Value = 4571.8945
$ws.Range("J31").Value = 4842.0303
$ws.Range("L31").Value = 4842.0303
$ws.Range("N31").Value = -5432.0303
# Row 34
$ws.Range("H34").Value = 4571.8945
$ws.Range("J34").Value = 4842.0303
$ws.Range("L34").Value = 4842.0303
$ws.Range("N34").Value = -5246.0303
# Row 107
$ws.Range("H107").Value = 37828.89
$ws.Range("I107").Value = 59360.707
$ws.Range("K107").Value = 59360.707
$ws.Range("M107").Value = -57440.707
# Row 122
$ws.Range("H122").Value = 3849388.5
$ws.Range("I122").Value = 5884951
$ws.Range("J122").Value = 4437.222
$ws.Range("K122").Value = 17654853
$ws.Range("L122").Value = 13311.666
$ws.Range("M122").Value = -17652403
$ws.Range("N122").Value = -18211.666
# Row 132
$ws.Range("H132").Value = 2966.4822
$ws.Range("I132").Value = 2875.2666
$ws.Range("J132").Value = 3339.6365
$ws.Range("K132").Value = 8625.799800000001
$ws.Range("L132").Value = 10018.9095
$ws.Range("M132").Value = -6095.799800000001
$ws.Range("N132").Value = -15078.9095

# ==== CUL sheet ====
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 2905.8647
$ws.Range("I68").Value = 1898.6
$ws.Range("J68").Value = 3278.926
$ws.Range("K68").Value = 5695.799999999999
$ws.Range("L68").Value = 9836.778
$ws.Range("M68").Value = -4884.799999999999
$ws.Range("N68").Value = -11458.778
# Row 71
$ws.Range("H71").Value = 2905.8647
$ws.Range("I71").Value = 1898.6
$ws.Range("J71").Value = 3278.926
$ws.Range("K71").Value = 17087.4
$ws.Range("L71").Value = 29510.334
$ws.Range("M71").Value = -13031.4
$ws.Range("N71").Value = -37622.334
# Row 98
$ws.Range("H98").Value = 3167.1667
$ws.Range("J98").Value = 3000
$ws.Range("L98").Value = 9000
$ws.Range("N98").Value = -11996
# Row 107
$ws.Range("H107").Value = 1197.7609
$ws.Range("I107").Value = 982.0714
$ws.Range("K107").Value = 2946.2142
$ws.Range("M107").Value = -1026.2142
# Row 122
$ws.Range("H122").Value = 323.72
$ws.Range("J122").Value = 288.16666
$ws.Range("L122").Value = 2593.49994
$ws.Range("N122").Value = -7493.49994

# ==== GSM sheet ====
$ws = $wb.Worksheets.Item("GSM")
# Row 39
$ws.Range("H39").Value = 68959
$ws.Range("J39").Value = 68959
$ws.Range("L39").Value = 68959
$ws.Range("N39").Value = -70023
# Row 70
$ws.Range("H70").Value = 14149.3
$ws.Range("I70").Value = 31016.363
$ws.Range("K70").Value = 31016.363
$ws.Range("M70").Value = -30746.363
# Row 73
$ws.Range("H73").Value = 14149.3
$ws.Range("I73").Value = 31016.363
$ws.Range("K73").Value = 31016.363
$ws.Range("M73").Value = -30080.363
# Row 102
$ws.Range("H102").Value = 2139.2188
$ws.Range("I102").Value = 2009.4615
$ws.Range("K102").Value = 2009.4615
$ws.Range("M102").Value = -387.4614999999999
# Row 113
$ws.Range("H113").Value = 69090.836
$ws.Range("I113").Value = 8777
$ws.Range("J113").Value = 99247.75
$ws.Range("K113").Value = 8777
$ws.Range("L113").Value = 99247.75
$ws.Range("M113").Value = -6607
$ws.Range("N113").Value = -103587.75
# Row 120
$ws.Range("H120").Value = 103342.164
$ws.Range("J120").Value = 103342.164
$ws.Range("L120").Value = 103342.164
$ws.Range("N120").Value = -113018.164
# Row 132
$ws.Range("H132").Value = 2362.9375
$ws.Range("I132").Value = 2111.162
$ws.Range("K132").Value = 6333.485999999999
$ws.Range("M132").Value = -3803.485999999999

# ==== LTW sheet ====
$ws = $wb.Worksheets.Item("LTW")
# Row 121
$ws.Range("H121").Value = 79271.60000000001
$ws.Range("J121").Value = 79271.60000000001
$ws.Range("L121").Value = 79271.60000000001
$ws.Range("N121").Value = -82765.60000000001
# Row 122
$ws.Range("H122").Value = 16416.5
$ws.Range("I122").Value = 18344.889
$ws.Range("J122").Value = 10631.333
$ws.Range("K122").Value = 55034.667
$ws.Range("L122").Value = 31893.999
$ws.Range("M122").Value = -52584.667
$ws.Range("N122").Value = -36793.999
# Row 132
$ws.Range("H132").Value = 172700.97
$ws.Range("I132").Value = 306035.34
$ws.Range("K132").Value = 918106.02
$ws.Range("M132").Value = -915576.02
# Row 136
$ws.Range("H136").Value = 4856.5293
$ws.Range("I136").Value = 4777.84
$ws.Range("K136").Value = 14333.52
$ws.Range("M136").Value = -11783.52
# Row 141
$ws.Range("H141").Value = 366936.75
$ws.Range("J141").Value = 366936.75
$ws.Range("L141").Value = 366936.75
$ws.Range("N141").Value = -377296.75

# ==== WVR sheet ====
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 70990.53
$ws.Range("I81").Value = 87071.5
$ws.Range("J81").Value = 6666.6665
$ws.Range("K81").Value = 174143
$ws.Range("L81").Value = 13333.333
$ws.Range("M81").Value = -173082
$ws.Range("N81").Value = -15455.333
# Row 84
$ws.Range("H84").Value = 70990.53
$ws.Range("I84").Value = 87071.5
$ws.Range("J84").Value = 6666.6665
$ws.Range("K84").Value = 870715
$ws.Range("L84").Value = 66666.66500000001
$ws.Range("M84").Value = -865411
$ws.Range("N84").Value = -77274.66500000001
# Row 96
$ws.Range("H96").Value = 4803
$ws.Range("I96").Value = 4766
$ws.Range("K96").Value = 4766
$ws.Range("M96").Value = -3393
# Row 97
$ws.Range("H97").Value = 54990
$ws.Range("J97").Value = 54990
$ws.Range("L97").Value = 54990
$ws.Range("N97").Value = -56972
# Row 121
$ws.Range("H121").Value = 51146.6
$ws.Range("J121").Value = 51146.6
$ws.Range("L121").Value = 51146.6
$ws.Range("N121").Value = -54640.6
# Row 132
$ws.Range("H132").Value = 16614.4
$ws.Range("I132").Value = 17650.814
$ws.Range("K132").Value = 52952.442
$ws.Range("M132").Value = -50422.442
